$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206, shifting existing rows 206:244 down to 207:245
$ws.Rows.Item(206).Insert()

# Populate the newly inserted row 206 with the new record's data
$ws.Range("A206").Value = 3
$ws.Range("B206").Value = "Femacal de La Calera"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 44476
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 100112017
$ws.Range("G206").Value = "Apio"
$ws.Range("H206").Value = "Americana (o)"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 60
$ws.Range("K206").Value = 9000
$ws.Range("L206").Value = 9000
$ws.Range("M206").Value = 9000
$ws.Range("N206").Value = "$/docena de matas"
$ws.Range("O206").Value = "Pan de Az$([char]0x00FA)car"
$ws.Range("P206").Value = 1500
$ws.Range("Q206").Value = 6
$ws.Range("R206").Value = "Hortaliza"
